$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # column G
    $val = $cell.Value2
    if ($null -ne $val -and $val -ne "") {
        $parts = $val -split ", "

        $hasSystem = $false
        foreach ($p in $parts) {
            if ($p.Equals("System")) { $hasSystem = $true }
        }

        $isLast = $false
        if ($parts.Count -gt 0) {
            $isLast = $parts[$parts.Count - 1].Equals("System")
        }

        if ($hasSystem -and -not $isLast) {
            $newParts = @()
            foreach ($p in $parts) {
                if (-not $p.Equals("System")) {
                    $newParts += $p
                }
            }
            $newParts += "System"
            $cell.Value = ($newParts -join ", ")
        }
    }
}
